$d = $word.ActiveDocument

# Locate the "Ratings:" table under the "PV Data:" section. It is the
# 3-row x 7-column table whose first cell reads "Conn" (Conn / Ph-1 / Ph-1 /
# Ph-2 / Ph-2 / Ph-3 / Ph-3 header, then kVA/pf sub-header, then the "D"
# data row ending in ...120/0.95/120/0.95/120/0.95).
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $firstCell = $tbl.Cell(1,1).Range.Text -replace "[\x07\x0d]", ""
    if ($firstCell -eq "Conn" -and $tbl.Rows.Count -eq 3 -and $tbl.Columns.Count -eq 7) {
        $targetTable = $tbl
    }
}

if ($targetTable -eq $null) {
    throw "Could not locate the PV Data Ratings table"
}

$insPoint = $d.Range($targetTable.Range.End, $targetTable.Range.End)

$newXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Dynamic Parameters</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:tblPr>
        <w:tblW w:w="0" w:type="auto"/>
        <w:jc w:val="center"/>
        <w:tblBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tblBorders>
        <w:tblLayout w:type="fixed"/>
        <w:tblCellMar>
          <w:left w:w="30" w:type="dxa"/>
          <w:right w:w="30" w:type="dxa"/>
        </w:tblCellMar>
        <w:tblLook w:val="0000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="662"/>
        <w:gridCol w:w="629"/>
        <w:gridCol w:w="629"/>
        <w:gridCol w:w="685"/>
        <w:gridCol w:w="592"/>
        <w:gridCol w:w="816"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="250"/>
          <w:jc w:val="center"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="662" w:type="dxa"/>
            <w:shd w:val="pct12" w:color="auto" w:fill="FFFFFF"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:lastRenderedPageBreak/>
              <w:t>R</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="629" w:type="dxa"/>
            <w:shd w:val="pct12" w:color="auto" w:fill="FFFFFF"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>X</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="629" w:type="dxa"/>
            <w:shd w:val="pct12" w:color="auto" w:fill="FFFFFF"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>Kp</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="685" w:type="dxa"/>
            <w:shd w:val="pct12" w:color="auto" w:fill="FFFFFF"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>kVDC</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="592" w:type="dxa"/>
            <w:shd w:val="pct12" w:color="auto" w:fill="FFFFFF"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>KP Tol</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="816" w:type="dxa"/>
            <w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>Safe Voltage</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="250"/>
          <w:jc w:val="center"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="662" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>0.5</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="629" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>0.5</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="629" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>0.01</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="685" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>0.03</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="592" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>0.1</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="816" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:snapToGrid w:val="0"/>
                <w:color w:val="000000"/>
                <w:sz w:val="16"/>
              </w:rPr>
              <w:t>0</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
'@

$insPoint.InsertXML($newXml)
